# Apply weekly price update to "Fruta, Vega Modelo de Temuco - Frambuesa" sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44574
$ws.Range("M2").Value = 200
$ws.Range("N2").Value = 3000
$ws.Range("O2").Value = 3000
$ws.Range("P2").Value = 3000
$ws.Range("S2").Value = 3000

# Row 3
$ws.Range("D3").Value = 44592
$ws.Range("M3").Value = 5
$ws.Range("N3").Value = 7500
$ws.Range("O3").Value = 7500
$ws.Range("P3").Value = 7500
$ws.Range("S3").Value = 7500

# Row 4
$ws.Range("D4").Value = 44999
$ws.Range("M4").Value = 25
$ws.Range("N4").Value = 2500
$ws.Range("O4").Value = 2500
$ws.Range("P4").Value = 2500
$ws.Range("S4").Value = 2500

# Row 5
$ws.Range("D5").Value = 44214
$ws.Range("M5").Value = 50
$ws.Range("N5").Value = 1800
$ws.Range("O5").Value = 1800
$ws.Range("P5").Value = 1800
$ws.Range("S5").Value = 1800

# Row 6
$ws.Range("D6").Value = 44323
$ws.Range("M6").Value = 20
$ws.Range("N6").Value = 3200
$ws.Range("O6").Value = 3200
$ws.Range("P6").Value = 3200
$ws.Range("R6").Value = 'Región de La Araucanía'
$ws.Range("S6").Value = 3200

# Row 7
$ws.Range("D7").Value = 44175
$ws.Range("M7").Value = 40
$ws.Range("N7").Value = 5000
$ws.Range("O7").Value = 5000
$ws.Range("P7").Value = 5000
$ws.Range("R7").Value = 'Provincia de Curicó'
$ws.Range("S7").Value = 5000

# Row 8
$ws.Range("D8").Value = 44616
$ws.Range("M8").Value = 200
$ws.Range("N8").Value = 3200
$ws.Range("O8").Value = 3200
$ws.Range("P8").Value = 3200
$ws.Range("R8").Value = 'Región de La Araucanía'
$ws.Range("S8").Value = 3200

# Row 9
$ws.Range("D9").Value = 44551
$ws.Range("M9").Value = 120
$ws.Range("N9").Value = 4500
$ws.Range("O9").Value = 4500
$ws.Range("P9").Value = 4500
$ws.Range("R9").Value = 'Región de O''Higgins'
$ws.Range("S9").Value = 4500

# Row 10
$ws.Range("D10").Value = 44998
$ws.Range("M10").Value = 20
$ws.Range("N10").Value = 2500
$ws.Range("O10").Value = 2500
$ws.Range("P10").Value = 2500
$ws.Range("S10").Value = 2500

# Row 11
$ws.Range("D11").Value = 44176
$ws.Range("N11").Value = 3000
$ws.Range("O11").Value = 3000
$ws.Range("P11").Value = 3000
$ws.Range("R11").Value = 'Región de O''Higgins'
$ws.Range("S11").Value = 3000

# Row 12
$ws.Range("D12").Value = 44567
$ws.Range("M12").Value = 80
$ws.Range("N12").Value = 2400
$ws.Range("O12").Value = 2400
$ws.Range("P12").Value = 2400
$ws.Range("S12").Value = 2400

# New row 13 (inserted data point)
$ws.Range("A13").Value = 10
$ws.Range("B13").Value = 'Vega Modelo de Temuco'
$ws.Range("C13").Value = 'La Araucanía'
$ws.Range("D13").Value = 44215
$ws.Range("E13").Value = 9
$ws.Range("F13").Value = 'Fruta'
$ws.Range("G13").Value = 100101
$ws.Range("H13").Value = 'Berries'
$ws.Range("I13").Value = 100101004
$ws.Range("J13").Value = 'Frambuesa'
$ws.Range("K13").Value = 'Sin especificar'
$ws.Range("L13").Value = 'Primera'
$ws.Range("M13").Value = 65
$ws.Range("N13").Value = 2800
$ws.Range("O13").Value = 2800
$ws.Range("P13").Value = 2800
$ws.Range("Q13").Value = '$/envase 1 kilo'
$ws.Range("R13").Value = 'Región de La Araucanía'
$ws.Range("S13").Value = 2800
$ws.Range("T13").Value = 1
$ws.Range("D13").NumberFormat = "YYYY-MM-DD HH:MM:SS"

